$d = $word.ActiveDocument

# --- Change 1: "A jar file can be found in the " + "dist" + " folder." -> one run
# (the original run is split across 3 runs with spell-check proofErr markers around
# "dist"; replacing the full sentence in one Find/Replace merges them into a single
# run and drops the proofErr markers.)
$d.Content.Find.Execute("A jar file can be found in the dist folder.", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    "A jar file can be found in the dist folder.", 2)

# --- Change 2: merge the run that only contains <w:br/> together with the
# proofErr-wrapped "Admin" run and the " password: pdc2023" run that follow it,
# producing a single run "<w:br/>Admin password: pdc2023", while leaving the
# preceding "Admin username: admin" run untouched.
$usernameRange = $d.Content
$usernameRange.Find.Execute("Admin username: admin", $true, $false, $false, `
    $false, $false, $true, 1, $false)
$usernameStart = $usernameRange.Start
$usernameEnd = $usernameRange.End

$tailRange = $d.Range($usernameEnd, $d.Content.End)
$tailRange.Find.Execute("Admin password: pdc2023", $true, $false, $false, `
    $false, $false, $true, 1, $false)
$tailEnd = $tailRange.End

# The merge range starts right at the <w:br/> run (immediately after "admin")
# and ends at the end of "pdc2023". Re-assigning FormattedText over runs that
# share identical formatting coalesces them into one run and clears any
# leftover spell-check proofErr markers.
$mergeRange = $d.Range($usernameEnd, $tailEnd)
$mergeRange.FormattedText = $mergeRange.FormattedText

# The coalesce above also pulls in the preceding "Admin username: admin" run
# (since it has identical formatting). Toggling a character property on just
# that leading run forces Word to split it back out into its own run, leaving
# the break+"Admin password: pdc2023" text merged together as desired.
$usernamePart = $d.Range($usernameStart, $usernameEnd)
$usernamePart.Bold = 1
$usernamePart.Bold = 0

# --- Change 3: append an empty paragraph followed by a "Github link: <url>"
# paragraph at the end of the document body.
$endPos = $d.Content.End
$tailAnchor = $d.Range($endPos - 1, $endPos - 1)
# Use a throw-away placeholder character "X" for the new blank paragraph so the
# paragraph-insert machinery doesn't need to leave a leftover empty run behind;
# the placeholder is deleted immediately afterwards, leaving a truly empty
# paragraph (matching how the rest of the document's blank paragraphs look).
$tailAnchor.InsertAfter("`rX`rGithub link: ")

$placeholder = $d.Content
$placeholder.Find.Execute("X", $true, $false, $false, $false, $false, $true, `
    1, $false)
$placeholderRange = $d.Range($placeholder.Start, $placeholder.End)
$placeholderRange.Delete()

$linkEndPos = $d.Content.End
$linkAnchor = $d.Range($linkEndPos - 1, $linkEndPos - 1)
$linkAnchor.InsertAfter("https://github.com/lahndrick/Shopping-System")

# Split "Github link: " and the URL into two separate runs (matching the
# target markup) by toggling a character property on just the label text.
$labelRange = $d.Content
$labelRange.Find.Execute("Github link: ", $true, $false, $false, $false, `
    $false, $true, 1, $false)
$label = $d.Range($labelRange.Start, $labelRange.End)
$label.Bold = 1
$label.Bold = 0
